$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert counter values (stop/restart/reset button logic reverted)
$ws.Range("B2").Value = 2
$ws.Range("B4").Value = 3

# Revert the barcode text in A4 back to the earlier (duplicate) value
$ws.Range("A4").Value = "a02403964220"

# Revert cell formatting on the barcode column: drop the custom "Malgun Gothic"
# font override and use left-aligned default font instead
$rng = $ws.Range("A2:A4")
$rng.ClearFormats()
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4108

# Restore selection
$ws.Range("B7").Select() | Out-Null
